# JS-Frameworks-Self-Evaluation-Protocol_tpPetkov.xlsx
# Update the "Issue page" score (row 25) from 20 to 15; the Total Score
# formula in C32 (=SUM(C6:C31)) recalculates automatically from 268 to 263.
# Also move the active selection from C19 to E25, matching the saved
# sheet view state in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C25: "Issue page" score, 20 -> 15
$ws.Range("C25").Value = 15

# Update the saved selection/active cell for the sheet view
$ws.Range("E25").Select()
